$d = $word.ActiveDocument

$replacements = @(
    @("2025-06-21 Saturday", "2025-06-22 Sunday"),
    @("991÷2=495, 1", "817÷5=163, 2"),
    @("884÷5=176, 4", "772÷5=154, 2"),
    @("334÷9=37, 1", "353÷6=58, 5"),
    @("952÷6=158, 4", "303÷9=33, 6"),
    @("948÷5=189, 3", "530÷3=176, 2"),
    @("397÷4=99, 1", "854÷3=284, 2"),
    @("489÷8=61, 1", "576÷8=72, 0"),
    @("317÷7=45, 2", "584÷2=292, 0"),
    @("210÷3=70, 0", "975÷4=243, 3"),
    @("319÷2=159, 1", "312÷8=39, 0"),
    @("314÷2=157, 0", "368÷7=52, 4"),
    @("584÷5=116, 4", "987÷3=329, 0"),
    @("998÷4=249, 2", "645÷8=80, 5"),
    @("504÷9=56, 0", "235÷8=29, 3"),
    @("778÷4=194, 2", "679÷4=169, 3"),
    @("988÷2=494, 0", "914÷9=101, 5"),
    @("174÷8=21, 6", "705÷5=141, 0"),
    @("247÷9=27, 4", "196÷2=98, 0"),
    @("949÷4=237, 1", "297÷7=42, 3"),
    @("138÷6=23, 0", "602÷4=150, 2"),
    @("608÷2=304, 0", "921÷9=102, 3"),
    @("821÷2=410, 1", "348÷9=38, 6"),
    @("996÷5=199, 1", "462÷2=231, 0"),
    @("214÷2=107, 0", "204÷5=40, 4"),
    @("839÷2=419, 1", "245÷2=122, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done applying $($replacements.Count) replacements"
